# Apply "Updated CAD models" pose data edit to Poses1.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update pose coordinate values (A2:C8)
$ws.Range("A2").Value = -2.50664
$ws.Range("B2").Value = 5.8799
$ws.Range("C2").Value = 7.86653

$ws.Range("A3").Value = 4.47819
$ws.Range("B3").Value = -8.9185
$ws.Range("C3").Value = 8.62735

$ws.Range("A4").Value = -2.2709
$ws.Range("B4").Value = 0.395378
$ws.Range("C4").Value = -4.81918

$ws.Range("A5").Value = -1.32134
$ws.Range("B5").Value = -7.50649
$ws.Range("C5").Value = -3.55002

$ws.Range("A6").Value = -3.5464
$ws.Range("B6").Value = -9.68343
$ws.Range("C6").Value = 9.79838

$ws.Range("A7").Value = 1.94974
$ws.Range("B7").Value = 3.61731
$ws.Range("C7").Value = 2.83973

$ws.Range("A8").Value = -2.83135
$ws.Range("B8").Value = -0.416536
$ws.Range("C8").Value = 6.56425

# Update the active selection shown on the sheet (A2:C9 range, anchored at A2)
$excel.Goto($ws.Range("A2:C9"))
